$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty row 31 (A31:C31) with the newly calculated
# payment-history values for this date.
$ws.Range("A31").Value = 44237
$ws.Range("B31").Value = 1810.0135264944515
$ws.Range("C31").Value = 225.54647350554853

# Update the active selection to reflect where the user ended up (E18).
$ws.Range("E18").Select()
